$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D, shifting existing D:K data to E:L.
$ws.Columns("D:D").Insert()

# The freshly inserted column has no number formatting of its own yet;
# copy the formats from column E (the old column D, now shifted right)
# so the new column D keeps using the same shared styles (date / number).
# Limit the copy to the sheet's used rows so we don't touch all 1M+ rows.
$ws.Range("E1:E102").Copy()
$ws.Range("D1:D102").PasteSpecial(-4122)  # xlPasteFormats

# New values for the inserted column D, keyed by row number.
$newD = @{
    7   = 43465
    8   = 142000
    13  = 0
    14  = 0
    15  = 0
    17  = 23000
    18  = 119000
    20  = -67200
    21  = 61600
    22  = 0
    23  = 51800
    24  = 9800
    25  = 0
    26  = 42000
    27  = 42000
    28  = 0
    30  = 0
    31  = 0
    32  = 67200
    33  = 42000
    34  = 0
    35  = 42000
    38  = 43465
    41  = 64900
    42  = 26500
    43  = 0
    44  = 0
    45  = 0
    46  = 0
    47  = 0
    48  = 48300
    49  = 55000
    50  = 0
    51  = 0
    53  = 0
    54  = 3363900
    57  = 25200
    58  = 0
    59  = 0
    60  = 0
    61  = 46200
    62  = 0
    63  = 0
    64  = 0
    65  = 0
    66  = 2988700
    68  = 0
    69  = 0
    70  = 0
    71  = 0
    72  = 75500
    73  = 0
    74  = 0
    75  = 0
    76  = 375200
    77  = 0
    80  = 43465
    81  = 42000
    83  = 9800
    84  = 0
    85  = 0
    86  = 0
    87  = 0
    88  = 0
    89  = 61700
    91  = -6300
    92  = 0
    93  = 0
    94  = -210500
    96  = -27500
    97  = 0
    98  = 0
    99  = 0
    100 = 24000
    101 = 0
    102 = -124700
}

foreach ($row in $newD.Keys) {
    $ws.Cells.Item($row, 4).Value = $newD[$row]
}

# Rows where the new column D should hold "NA" text instead of a number.
$newDNA = @(9, 10, 12, 29, 52)
foreach ($row in $newDNA) {
    $ws.Cells.Item($row, 4).Value = "NA"
}
